$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "23.098.83"
$ws.Range("E2").Value = "  -3.42%  "
$ws.Range("D3").Value = "1.601.64"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'301.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("D7").Value = "'0.3783"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("D8").Value = "'0.3654"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.28%  "
$ws.Range("D9").Value = "'49.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("D10").Value = "'1.271"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.00%  "
$ws.Range("D11").Value = "'0.08159"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  -3.51%  "
$ws.Range("D14").Value = "'6.600"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.86%  "
$ws.Range("D15").Value = "'0.00001261"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").Value = "'7.401"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.42%  "
$ws.Range("D17").Value = "1.604.35"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "'91.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").Value = "'0.06837"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "'18.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.67%  "
$ws.Range("D21").Value = "'6.581"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("D22").Value = "'0.5570"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.38%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'13.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.57%  "
$ws.Range("D25").Value = "23.114.35"
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("D27").Value = "'2.722"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.84%  "
$ws.Range("D28").Value = "'21.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Value = "'5.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "'132.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("D32").Value = "'2.397"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").Value = "'6.850"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.04%  "
$ws.Range("D34").Value = "1.781.19"
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("D35").Value = "'0.9617"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.97%  "
$ws.Range("D36").Value = "'0.07696"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.85%  "
$ws.Range("D37").Value = "'6.294"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.48%  "
$ws.Range("D38").Value = "'0.02734"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.51%  "
$ws.Range("D39").Value = "'0.2557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.69%  "
$ws.Range("D40").Value = "'0.08914"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("D41").Value = "'10.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.63%  "
$ws.Range("D42").Value = "'1.371"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("D43").Value = "'0.7107"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.30%  "
$ws.Range("E44").Value = "  -6.44%  "
$ws.Range("D45").Value = "'15.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.97%  "
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("D47").Value = "'2.315"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.08%  "
$ws.Range("D49").Value = "'3.991"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("D50").Value = "'132.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("E51").Value = "  -4.36%  "
